$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.605.61"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "3.915.23"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'603.71"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'165.51"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").Value = "3.913.85"
$ws.Range("E7").Value = "  +3.09%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("D11").Value = "'6.40"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "'37.30"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "4.564.39"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "3.910.23"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "68.722.07"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "'7.44"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'17.23"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "'10.99"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").Value = "'485.95"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "'0.0000169"
$ws.Range("E24").Value = "  +11.90%  "
$ws.Range("D25").Value = "'84.68"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "'12.09"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'2.93"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").Value = "4.066.01"
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("D32").Value = "'2.39"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").Value = "'7.77"
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("D34").Value = "'32.00"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "3.862.05"
$ws.Range("E35").Value = "  +3.05%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").Value = "'5.93"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.16"
$ws.Range("E40").Value = "  +4.46%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'0.316"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").Value = "'428.79"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("D44").Value = "'48.39"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'8.50"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'26.22"
$ws.Range("E48").Value = "  +7.21%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'141.85"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.813.39"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0353"
$ws.Range("E51").Value = "  +0.72%  "
